# Insert a new data row at row 376 (pushing the existing rows 376-433 down
# to 377-434) and populate it with the new weekly record, matching the
# author's diff for "Hortaliza, Vega Modelo de Temuco - Ciboulette".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift everything at/after row 376 down by one row.
$ws.Rows.Item(376).Insert()

# Populate the newly inserted row with the new record's values.
$ws.Cells.Item(376, 1).Value  = 10
$ws.Cells.Item(376, 2).Value  = "Vega Modelo de Temuco"
$ws.Cells.Item(376, 3).Value  = "La Araucanía"
$ws.Cells.Item(376, 4).Value  = 45218
$ws.Cells.Item(376, 5).Value  = 9
$ws.Cells.Item(376, 6).Value  = 100112039
$ws.Cells.Item(376, 7).Value  = "Ciboulette"
$ws.Cells.Item(376, 8).Value  = "Sin especificar"
$ws.Cells.Item(376, 9).Value  = "Primera"
$ws.Cells.Item(376, 10).Value = 85
$ws.Cells.Item(376, 11).Value = 7000
$ws.Cells.Item(376, 12).Value = 7000
$ws.Cells.Item(376, 13).Value = 7000
$ws.Cells.Item(376, 14).Value = "$/docena de atados"
$ws.Cells.Item(376, 15).Value = "Provincia de Cautín"
$ws.Cells.Item(376, 16).Value = 2333
$ws.Cells.Item(376, 17).Value = 3
$ws.Cells.Item(376, 18).Value = "Hortaliza"
